$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Code")

$rows = @{
  3 = @("Lars Johansson", "O-XFT", "Dsgnr1")
  4 = @("Anna Lanner", "PjM", "O2")
  5 = @("Lena Serra", "CCB", "CCB-A1")
  6 = @("Marja Hammarstrand", "TC", "TC1")
  7 = @("Magnus Bergh", "OPO", "OPO2")
  8 = @("Björn Ternström", "APO", "APO1")
  9 = @("Elizabeth Hansson", "KM", "CAKM1")
  10 = @("Marcel van Torgeren", "PgM", "PgM2")
  11 = @("Helena Eberil", "SM", "O6")
  12 = @("Simeon Soetan", "O-XFT", "Dsgnr2")
  13 = @("Marcus Larsson", "OPO", "OPO3")
  14 = @("Henric Stenhoff", "DM", "DM2")
  15 = @("Patrik Främme", "SM", "SM2")
  16 = @("Björn Östlund", "TM", "TM1")
  17 = @("Lars Rundberg", "O-XFT", "Dsgnr3")
  18 = @("Eva Cullman", "SM", "O7")
  19 = @("Mikael Krekola", "PG", "PG3")
  20 = @("Pierre Svärd", "O-XFT", "Dsgnr4")
  21 = @("Anny Lei", "FjPM", "FPjL1")
  22 = @("Jun Johansson", "RBS", "PdMRBS1")
  23 = @("Per Simonsson", "RBS", "RBS1")
  24 = @("Henrik Sundh", "SPM", "SPM1")
  25 = @("Jeanette Munro", "O-XFT", "Dsgnr5")
  26 = @("Ricardo Morales", "O-XFT", "Dsgnr6")
  27 = @("Carl Ohlsson", "O-XFT", "Dsgnr7")
  28 = @("Niklas Isaksson", "PgM", "PgM3")
  29 = @("Thomas Andersson", "SrM", "SrM1")
  30 = @("Per Lofter", "SM", "SM3")
  31 = @("Thomas Nyberg", "OPO", "OPO4")
  32 = @("Lena Doverfors", "CCB", "CCB-A2")
  33 = @("Ulf Olsson", "CPI", "CPIW1")
  34 = @("Stefan Jigsved", "SM", "SM4")
  35 = @("Fredrik Svanfeldt", "DM", "DM3")
  36 = @("Sven-Eric Ericson", "O-XFT", "Dsgnr8")
  37 = @("Mats Eriksson", "CL", "CL1")
  38 = @("Eva Telandersson", "TR", "TR-A1")
  39 = @("Anna Ekedahl", "CPI", "CPIPjL1")
  40 = @("Juhan Zhao", "CPI", "CPI1")
  41 = @("Niclas Fremling", "O-XFT", "Dsgnr9")
  42 = @("Mats Nyrenius", "O-XFT", "Dsgnr10")
  43 = @("Gunnar Lindh", "CPI", "CPIW2")
  44 = @("Sari Eklund", "CPI", "CPIW3")
  45 = @("Irina Romanova", "RBS", "RBS2")
  46 = @("Karin Åkesson", "CIRV", "Cirv1")
  47 = @("Peter Astrof", "O-XFT", "Dsgnr11")
  48 = @("Malin Aguilera", "O-XFT", "Dsgnr12")
  49 = @("Henrik Larsson", "O-XFT", "Dsgnr13")
  50 = @("Thomas Stephanson", "ITS", "ITS1")
  51 = @("Suxia Oldemark", "PG", "PG4")
  52 = @("Anders Borghed", "PG", "PG5")
  53 = @("Peter Malmqvist", "O-XFT", "Dsgnr14")
  54 = @("Malin Lind", "O-XFT", "Dsgnr15")
  55 = @("Bengt Skarin", "O-XFT", "Dsgnr16")
  56 = @("Michael Nordquist", "O-XFT", "Dsgnr17")
}

foreach ($r in $rows.Keys) {
  $vals = $rows[$r]
  $ws.Cells.Item([int]$r, 6).Value = $vals[0]
  $ws.Cells.Item([int]$r, 7).Value = $vals[1]
  $ws.Cells.Item([int]$r, 8).Value = $vals[2]
}

$ws.Range("H47").Select()

Write-Output "done"
